$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5286.8
$ws.Range("I86").Value = 1882
$ws.Range("J86").Value = 10394
$ws.Range("K86").Value = 1882
$ws.Range("L86").Value = 10394
$ws.Range("M86").Value = -759
$ws.Range("N86").Value = -12640
$ws.Range("H89").Value = 5286.8
$ws.Range("I89").Value = 1882
$ws.Range("J89").Value = 10394
$ws.Range("K89").Value = 9410
$ws.Range("L89").Value = 51970
$ws.Range("M89").Value = -3794
$ws.Range("N89").Value = -63202
$ws.Range("H100").Value = 1534.7
$ws.Range("J100").Value = 2127.4
$ws.Range("L100").Value = 2127.4
$ws.Range("N100").Value = -3209.4
$ws.Range("H111").Value = 11606.714
$ws.Range("I111").Value = 32014.5
$ws.Range("J111").Value = 3443.6
$ws.Range("K111").Value = 96043.5
$ws.Range("L111").Value = 10330.8
$ws.Range("M111").Value = -92976.5
$ws.Range("N111").Value = -16464.8
$ws.Range("H125").Value = 3714.8
$ws.Range("I125").Value = 3908
$ws.Range("J125").Value = 3586
$ws.Range("K125").Value = 35172
$ws.Range("L125").Value = 32274
$ws.Range("M125").Value = -32712
$ws.Range("N125").Value = -37194
$ws.Range("H132").Value = 10423021
$ws.Range("I132").Value = 11370342
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 34111026
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -34108496
$ws.Range("N132").Value = -12560

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1297.6897
$ws.Range("I61").Value = 1212.2693
$ws.Range("J61").Value = 2038
$ws.Range("K61").Value = 1212.2693
$ws.Range("L61").Value = 2038
$ws.Range("M61").Value = -1000.2693
$ws.Range("N61").Value = -2462
$ws.Range("H64").Value = 38266.668
$ws.Range("J64").Value = 38266.668
$ws.Range("L64").Value = 38266.668
$ws.Range("N64").Value = -38762.668
$ws.Range("H67").Value = 38266.668
$ws.Range("J67").Value = 38266.668
$ws.Range("L67").Value = 38266.668
$ws.Range("N67").Value = -39982.668
$ws.Range("H74").Value = 2277.4285
$ws.Range("I74").Value = 1141.6666
$ws.Range("K74").Value = 1141.6666
$ws.Range("M74").Value = -267.6666
$ws.Range("H77").Value = 2277.4285
$ws.Range("I77").Value = 1141.6666
$ws.Range("K77").Value = 5708.333000000001
$ws.Range("M77").Value = -1340.333000000001
$ws.Range("H105").Value = 46240
$ws.Range("J105").Value = 46240
$ws.Range("L105").Value = 46240
$ws.Range("N105").Value = -53228
$ws.Range("H136").Value = 1297.6897
$ws.Range("I136").Value = 1212.2693
$ws.Range("J136").Value = 2038
$ws.Range("K136").Value = 3636.8079
$ws.Range("L136").Value = 6114
$ws.Range("M136").Value = -1086.8079
$ws.Range("N136").Value = -11214

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 36092.793
$ws.Range("I20").Value = 51729.65
$ws.Range("K20").Value = 51729.65
$ws.Range("M20").Value = -51482.65
$ws.Range("H62").Value = 50387
$ws.Range("J62").Value = 50387
$ws.Range("L62").Value = 50387
$ws.Range("N62").Value = -51759
$ws.Range("H65").Value = 50387
$ws.Range("J65").Value = 50387
$ws.Range("L65").Value = 151161
$ws.Range("N65").Value = -158025
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23303.666
$ws.Range("I31").Value = 47004.727
$ws.Range("J31").Value = 3248.923
$ws.Range("K31").Value = 47004.727
$ws.Range("L31").Value = 3248.923
$ws.Range("M31").Value = -46709.727
$ws.Range("N31").Value = -3838.923
$ws.Range("H34").Value = 23303.666
$ws.Range("I34").Value = 47004.727
$ws.Range("J34").Value = 3248.923
$ws.Range("K34").Value = 47004.727
$ws.Range("L34").Value = 3248.923
$ws.Range("M34").Value = -46802.727
$ws.Range("N34").Value = -3652.923
$ws.Range("H132").Value = 31917492
$ws.Range("I132").Value = 30305312
$ws.Range("K132").Value = 90915936
$ws.Range("M132").Value = -90913406

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 12418.223
$ws.Range("J5").Value = 15410.286
$ws.Range("L5").Value = 46230.858
$ws.Range("N5").Value = -46454.858
$ws.Range("H14").Value = 449.42856
$ws.Range("I14").Value = 449.42856
$ws.Range("K14").Value = 1348.28568
$ws.Range("M14").Value = -1175.28568
$ws.Range("H48").Value = 1000
$ws.Range("J48").Value = 1000
$ws.Range("L48").Value = 3000
$ws.Range("N48").Value = -3500
$ws.Range("H86").Value = 25001062
$ws.Range("I86").Value = 366.66666
$ws.Range("J86").Value = 40001480
$ws.Range("K86").Value = 1099.99998
$ws.Range("L86").Value = 120004440
$ws.Range("M86").Value = 86.00001999999995
$ws.Range("N86").Value = -120006812
$ws.Range("H89").Value = 25001062
$ws.Range("I89").Value = 366.66666
$ws.Range("J89").Value = 40001480
$ws.Range("K89").Value = 3299.99994
$ws.Range("L89").Value = 360013320
$ws.Range("M89").Value = 2628.00006
$ws.Range("N89").Value = -360025176
$ws.Range("H100").Value = 6667.5386
$ws.Range("J100").Value = 6667.5386
$ws.Range("L100").Value = 20002.6158
$ws.Range("N100").Value = -21624.6158
$ws.Range("H131").Value = 725.83
$ws.Range("I131").Value = 432.78946
$ws.Range("J131").Value = 794.5679
$ws.Range("K131").Value = 1298.36838
$ws.Range("L131").Value = 2383.7037
$ws.Range("M131").Value = 3741.63162
$ws.Range("N131").Value = -12463.7037
$ws.Range("H135").Value = 12418.223
$ws.Range("J135").Value = 15410.286
$ws.Range("L135").Value = 138692.574
$ws.Range("N135").Value = -143762.574

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 72340.60000000001
$ws.Range("I70").Value = 122716.94
$ws.Range("J70").Value = 6463.846
$ws.Range("K70").Value = 122716.94
$ws.Range("L70").Value = 6463.846
$ws.Range("M70").Value = -122446.94
$ws.Range("N70").Value = -7003.846
$ws.Range("H73").Value = 72340.60000000001
$ws.Range("I73").Value = 122716.94
$ws.Range("J73").Value = 6463.846
$ws.Range("K73").Value = 122716.94
$ws.Range("L73").Value = 6463.846
$ws.Range("M73").Value = -121780.94
$ws.Range("N73").Value = -8335.846
$ws.Range("H80").Value = 71430630
$ws.Range("I80").Value = 142859000
$ws.Range("J80").Value = 2247.1428
$ws.Range("K80").Value = 142859000
$ws.Range("L80").Value = 2247.1428
$ws.Range("M80").Value = -142858002
$ws.Range("N80").Value = -4243.1428
$ws.Range("H83").Value = 71430630
$ws.Range("I83").Value = 142859000
$ws.Range("J83").Value = 2247.1428
$ws.Range("K83").Value = 714295000
$ws.Range("L83").Value = 11235.714
$ws.Range("M83").Value = -714290008
$ws.Range("N83").Value = -21219.714

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 49252.24
$ws.Range("I40").Value = 167766.5
$ws.Range("J40").Value = 1846.5333
$ws.Range("K40").Value = 167766.5
$ws.Range("L40").Value = 1846.5333
$ws.Range("M40").Value = -167630.5
$ws.Range("N40").Value = -2118.5333
$ws.Range("H61").Value = 2460.6155
$ws.Range("I61").Value = 1572.8572
$ws.Range("J61").Value = 3496.3333
$ws.Range("K61").Value = 1572.8572
$ws.Range("L61").Value = 3496.3333
$ws.Range("M61").Value = -1370.8572
$ws.Range("N61").Value = -3900.3333
$ws.Range("H113").Value = 2460.6155
$ws.Range("I113").Value = 1572.8572
$ws.Range("J113").Value = 3496.3333
$ws.Range("K113").Value = 1572.8572
$ws.Range("L113").Value = 3496.3333
$ws.Range("M113").Value = 597.1428000000001
$ws.Range("N113").Value = -7836.3333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1154.1842
$ws.Range("I132").Value = 772.6129
$ws.Range("J132").Value = 2844
$ws.Range("K132").Value = 2317.8387
$ws.Range("L132").Value = 8532
$ws.Range("M132").Value = 212.1613000000002
$ws.Range("N132").Value = -13592
